$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# Extend formatting for the new rows (11-18) by copying the format of an
# existing fully-styled data row (row 10) down across the new range first,
# so every new cell keeps the same fill/border style (s="3") as the rest
# of the table.
$srcFormatRow = $ws.Range("A10:E10")
$destFormatRange = $ws.Range("A11:E18")
$srcFormatRow.Copy($destFormatRange)

# Target contents for rows 7-18 of the "Test Steps" sheet (columns A-E).
$rows = @{
    7  = @("searchListingPageTestCases", "", "waitForElementPresent", "search_box", "")
    8  = @("searchListingPageTestCases", "", "click", "search_box", "")
    9  = @("searchListingPageTestCases", "", "type", "search_box", "apple")
    10 = @("searchListingPageTestCases", "", "click", "", "apple")
    11 = @("searchListingPageTestCases", "", "addProduct", "search_list_projectNames|search_list_addBtn", "Apple - Royal Gala")
    12 = @("searchListingPageTestCases", "", "verifyText", "productName_text", "Apple - Royal Gala")
    13 = @("searchListingPageTestCases", "", "click", "product_add_btn", "")
    14 = @("searchListingPageTestCases", "", "goBack", "", "")
    15 = @("searchListingPageTestCases", "", "addProduct", "search_list_projectNames|search_list_addBtn", "Green Apple")
    16 = @("searchListingPageTestCases", "", "verifyText", "productName_text", "Green Apple")
    17 = @("searchListingPageTestCases", "", "click", "product_add_btn", "")
    18 = @("searchListingPageTestCases", "", "click", "search_list_footerpage", "")
}

# Write the rows in the same order the original author typed them in (rows
# 7-13 in sequence, then row 18, then the remaining 14-17) so that brand
# new values land in the shared-string table in the same order they do in
# the reference workbook.
$writeOrder = @(7, 8, 9, 10, 11, 12, 13, 18, 14, 15, 16, 17)
foreach ($r in $writeOrder) {
    $vals = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Match the author's final view state: scrolled so column C is the
# left-most visible column, with E12 selected.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E12").Select()
